$wb = $excel.ActiveWorkbook

# Escapement (sheet1)
$ws = $wb.Worksheets.Item("Escapement")
$ws.Range("A60").Value = 878968.868471029
$ws.Range("B60").Value = 4410596.35098081
$ws.Range("C60").Value = 4696645.96731127
$ws.Range("A61").Value = 378746.452214836
$ws.Range("B61").Value = 3747413.98088953
$ws.Range("C61").Value = 3455654.72002602

# Total Catch (sheet2)
$ws = $wb.Worksheets.Item("Total Catch")
$ws.Range("A42").Value = 298096.32803626
$ws.Range("B42").Value = 4773590.18384449
$ws.Range("C42").Value = 1767566.87033872
$ws.Range("A43").Value = 158705.560643679
$ws.Range("B43").Value = 4887082.97075061
$ws.Range("C43").Value = 1518614.82166643
$ws.Range("A47").Value = 834667.767286461
$ws.Range("B47").Value = 5080167.23898994
$ws.Range("C47").Value = 1416936.51016396
$ws.Range("A51").Value = 313863.86832332
$ws.Range("B51").Value = 1934907.65872441
$ws.Range("C51").Value = 605771.038351433
$ws.Range("A52").Value = 305585.210157119
$ws.Range("B52").Value = 2103520.26295737
$ws.Range("C52").Value = 1116060.72448932
$ws.Range("A53").Value = 1096126.91850604
$ws.Range("B53").Value = 4405280.75965445
$ws.Range("C53").Value = 892371.696285627
$ws.Range("A54").Value = 992153.824034883
$ws.Range("B54").Value = 3078111.22601249
$ws.Range("C54").Value = 1679176.19875729
$ws.Range("A60").Value = 1144216.44797035
$ws.Range("B60").Value = 8260160.36212417
$ws.Range("C60").Value = 10376489.3239696
$ws.Range("A61").Value = 2617030.38634555
$ws.Range("B61").Value = 8706863.49050612
$ws.Range("C61").Value = 12939697.2338306

# Run Size (sheet3)
$ws = $wb.Worksheets.Item("Run Size")
$ws.Range("A42").Value = 492184.328002102
$ws.Range("B42").Value = 6233372.18409829
$ws.Range("C42").Value = 2409659.77037572
$ws.Range("A43").Value = 268355.560684681
$ws.Range("B43").Value = 6430474.9706263
$ws.Range("C43").Value = 2062487.12170713
$ws.Range("A47").Value = 1889371.76746346
$ws.Range("B47").Value = 6804843.23876095
$ws.Range("C47").Value = 1909482.51020936
$ws.Range("A51").Value = 507189.868351323
$ws.Range("B51").Value = 2699118.65871007
$ws.Range("C51").Value = 1038209.03833244
$ws.Range("A52").Value = 692621.210130621
$ws.Range("B52").Value = 3286868.26306597
$ws.Range("C52").Value = 2010208.72445032
$ws.Range("A53").Value = 1436716.91855654
$ws.Range("B53").Value = 7169894.75947345
$ws.Range("C53").Value = 1510848.69638763
$ws.Range("A54").Value = 1643325.82406839
$ws.Range("B54").Value = 5019585.22647469
$ws.Range("C54").Value = 2475860.19868229
$ws.Range("A60").Value = 2023168.44794138
$ws.Range("B60").Value = 12670316.362105
$ws.Range("C60").Value = 15073788.3242808
$ws.Range("A61").Value = 2995798.38626038
$ws.Range("B61").Value = 12454475.4913957
$ws.Range("C61").Value = 16394969.2338566

# Run Size no Offshore (sheet4)
$ws = $wb.Worksheets.Item("Run Size no Offshore")
$ws.Range("A60").Value = 1943383.91943144
$ws.Range("B60").Value = 12039622.1453228
$ws.Range("C60").Value = 14289207.767524
$ws.Range("A61").Value = 2869435.5151566
$ws.Range("B61").Value = 11870098.4776175
$ws.Range("C61").Value = 15584696.7552589
